# process executive data in test2.py
# Update computed percentage values on the active sheet to reflect the
# reprocessed executive data (counts/denominators changed for several
# categories, shifting percentages across rows/columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 12.5
$ws.Range("B2").Value = 10.25641025641026

$ws.Range("A3").Value = 25
$ws.Range("B3").Value = 28.2051282051282

$ws.Range("A4").Value = 37.5
$ws.Range("B4").Value = 29.48717948717949

$ws.Range("A5").Value = 25
$ws.Range("B5").Value = 32.05128205128205
$ws.Range("C5").Value = 28.08988764044944
$ws.Range("D5").Value = 13.58024691358025

$ws.Range("C6").Value = 15.73033707865168
$ws.Range("D6").Value = 30.8641975308642

$ws.Range("C7").Value = 56.17977528089888
$ws.Range("D7").Value = 55.55555555555556

$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("E10").Value = 49.25373134328358
$ws.Range("F10").Value = 40

$ws.Range("E11").Value = 50.74626865671642
$ws.Range("F11").Value = 60
